# Weekly refresh: a new week's record is inserted at the top of the
# Brocoli data block (row 63) and all the subsequent rows shift down by
# one, extending the used range from A1:R163 to A1:R164.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure the new bottom row (164) gets the same date number format as
# the rest of column D before it receives a value, so Excel doesn't
# invent a brand-new style for it.
$ws.Range("D64:D164").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# 1) Capture the existing block (rows 63..163, all columns A..R) and
#    shift it down one row (into 64..164) in a single bulk operation so
#    every column (dates, numbers, text) moves together.
$srcRange = $ws.Range("A63:R163")
$values = $srcRange.Value()
$dstRange = $ws.Range("A64:R164")
$dstRange.Value = $values

# 2) Row 63 becomes the new week's record. It keeps the same
#    market/quality/price/origin data it already had (that data simply
#    moved down to row 64 as part of the bulk shift above); only its
#    date advances to the new reporting date.
$ws.Range("D63").Value = 44467
